$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.96555634561868, 50.07445654262091]"
$ws.Range("T2").Value = "[49.98107618869848, 50.05870093046546]"
$ws.Range("L3").Value = "[49.96166040260586, 50.10894124825985]"
$ws.Range("T3").Value = "[49.95446495664381, 50.03177752900917]"
